$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

function Delete-Text($old) {
    $r = $d.Content.Duplicate
    $r.Find.Execute($old, $true, $true, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
    $r.Delete()
}

# --- Title ---
Replace-Text "Unveiling the Encryption Enigma" `
    "Unveiling the Symphony of Mathematics: A Journey Through Patterns, Relationships, and Problem-Solving"

# --- Author name ---
Replace-Text "Richard Stallman" "Miss Alida Bustamante"

# --- Email address ---
# original runs: "rms@stallman" / "." / "org"
# new content:   "alida"        / "." / "bustamante@robomentor" / "." / "net"
Replace-Text "rms@stallman" "alida"
Replace-Text "org" "bustamante@robomentor.net"

# --- Intro paragraph sentences (separated by <w:br/> into 3 blocks) ---
Replace-Text "Cryptography, the art of concealing information, has played a pivotal role in safeguarding data since ancient times" `
    "Mathematics, the universal language of science, unravels the hidden order and patterns that govern our world"

Replace-Text " From the enigmatic hieroglyphs of ancient Egypt to the modern-day digital algorithms, encryption has evolved as a cornerstone of secure communication" `
    " It is a symphony of numbers, shapes, and relationships, beckoning us to explore the intricacies of the cosmos"

Replace-Text " In this essay, we delve into the depths of encryption, exploring its historical milestones, unraveling its underlying principles, and examining its far-reaching applications in various fields" `
    " From the intricate designs of nature to the mind-boggling complexities of quantum physics, mathematics provides the tools to decipher the enigmatic puzzles that surround us"

Replace-Text "The genesis of encryption can be traced back to civilizations that sought to protect confidential messages from prying eyes" `
    "As we embark on this mathematical odyssey, we will uncover the fundamental principles that underpin this fascinating discipline"

Replace-Text " Early techniques, such as the Caesar cipher, employed simple character substitutions, laying the foundation for more sophisticated methods" `
    " We will explore the world of numbers, learning about their properties and operations"

Replace-Text " With the advent of the digital age, encryption underwent a transformative shift" `
    " We will delve into the realm of geometry, discovering the beauty and symmetry of shapes and their relationships"

Replace-Text " The introduction of public-key cryptography, epitomized by the RSA algorithm, revolutionized secure communication, enabling the exchange of encrypted messages without the need for a shared secret key" `
    " Algebra will introduce us to the power of variables and equations, enabling us to model and solve real-world problems"

Replace-Text "Encryption finds widespread application across diverse domains, ranging from military and diplomatic communications to financial transactions and digital signatures" `
    "Statistics, the science of data analysis, will equip us with the skills to interpret and make sense of the vast amounts of information that permeate our modern world"

Replace-Text " Its significance extends to safeguarding sensitive information in healthcare records, protecting intellectual property, and ensuring the integrity of electronic voting systems" `
    " Calculus, the mathematics of change, will open up new avenues of understanding, revealing the intricate dance of motion and the interplay of forces"

Replace-Text " The advent of quantum computing poses new challenges to traditional encryption methods, necessitating the exploration of post-quantum algorithms to maintain the inviolability of encrypted data" `
    " Through these explorations, we will discover the elegance and power of mathematics, appreciating its role as a fundamental pillar of scientific inquiry and technological advancement"

# --- Summary paragraph ---
Replace-Text "Encryption has evolved from ancient ciphers to sophisticated digital algorithms, playing a vital role in protecting information in various spheres" `
    "Our mathematical journey has taken us through the captivating world of numbers, shapes, relationships, and problem-solving"

Replace-Text " Its historical milestones, encompassing the Caesar cipher and public-key cryptography, have shaped the landscape of secure communication" `
    " We have explored the fundamental principles of arithmetic, geometry, algebra, statistics, and calculus, gaining a deeper appreciation for the beauty and power of mathematics"

Replace-Text " Encryption's applications span military, finance, healthcare, intellectual property, and electronic voting, emphasizing its multifaceted importance" `
    " Along the way, we have discovered the diverse applications of mathematics in science, engineering, technology, and everyday life"

# Drop the remaining summary sentences (quantum-computing caveat + closing line),
# keeping the final period that ends the paragraph.
Delete-Text " While quantum computing poses challenges, the pursuit of post-quantum algorithms ensures the continued efficacy of encryption in the face of emerging threats. Encryption remains an indispensable tool for safeguarding sensitive data and upholding privacy in the digital realm"

# --- Add a trailing empty paragraph at the very end of the document body ---
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

Write-Output "done"
